$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.626.89"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "2.049.03"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.51"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.667"
$ws.Range("E6").Value = "  +1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "56.97"
$ws.Range("E7").Value = "  -2.23%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "63.20"
$ws.Range("E9").Value = "  +6.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.368"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  -3.55%  "

$ws.Range("E12").Value = "  -3.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.920"
$ws.Range("E13").Value = "  +4.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.66"
$ws.Range("E14").Value = "  -3.83%  "

$ws.Range("D15").Value = "2.349.38"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("D17").Value = "2.061.33"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.67"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").Value = "36.551.17"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.03"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -3.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.88"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("E23").Value = "  -4.93%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  +2.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  -6.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.83"
$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.02"
$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").Value = "  -2.08%  "

$ws.Range("E31").Value = "  +3.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.06"
$ws.Range("E32").Value = "  -7.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0603"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  -9.42%  "

$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0870"
$ws.Range("E36").Value = "  +1.72%  "

$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("E38").Value = "  -7.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.12"
$ws.Range("E39").Value = "  +4.03%  "

$ws.Range("E40").Value = "  -6.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -1.98%  "

$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("E43").Value = "  -4.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.03"
$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0910"
$ws.Range("E45").Value = "  -5.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.97"
$ws.Range("E46").Value = "  -4.47%  "

$ws.Range("D47").Value = "1.379.03"
$ws.Range("E47").Value = "  +4.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.46"
$ws.Range("E48").Value = "  +8.97%  "

$ws.Range("E49").Value = "  +2.65%  "

$ws.Range("E50").Value = "  -4.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.85"
$ws.Range("E51").Value = "  +1.06%  "
